# Fill in the "comp_amp variations (w/o FPU)" section of the Profiling
# sheet (rows 18-21), mirroring the existing SM / GeLU variation blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Profiling")

# Row 18 - LOG_AMP_FXP_LUT
$ws.Range("B18").Value = "'"
$ws.Range("C18").Value = "LOG_AMP_FXP_LUT"
$ws.Range("D18").Value = 767747770
$ws.Range("F18").Value = 92257
$ws.Range("G18").Value = 2229

# Row 19 - LOG_AMP_FXP_APPROX
$ws.Range("C19").Value = "LOG_AMP_FXP_APPROX"
$ws.Range("D19").Value = 728430738
$ws.Range("F19").Value = 215213
$ws.Range("G19").Value = 24505

# Row 20 - sub-header: "all fixed (although log_amp is faulty)"
$ws.Range("B20").Value = "all fixed (although log_amp is faulty)"
$ws.Range("C20").Value = "'"

# Row 21 - logamp_approx_sm_fixed_gelu_pwl
$ws.Range("C21").Value = "logamp_approx_sm_fixed_gelu_pwl"
$ws.Range("D21").Value = 294178509
$ws.Range("F21").Value = 208442
$ws.Range("G21").Value = 19958
